$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.45799994422025
$ws.Range("C2").Value = 9.160868140125007
$ws.Range("D2").Value = 7.747803205339658
$ws.Range("E2").Value = 9.868800480874389
$ws.Range("F2").Value = 40.75843688601618
$ws.Range("I2").Value = 32.27767050510513
$ws.Range("L2").Value = 10.3508525004283
$ws.Range("M2").Value = 16.68056971028254
$ws.Range("B3").Value = 18.1203113214724
$ws.Range("C3").Value = 8.549156298291301
$ws.Range("D3").Value = 7.761444272747956
$ws.Range("E3").Value = 9.853265667986904
$ws.Range("F3").Value = 40.40164030969947
$ws.Range("I3").Value = 32.17800951951017
$ws.Range("L3").Value = 10.36146805487886
$ws.Range("M3").Value = 16.63224542646559
$ws.Range("B4").Value = 17.91599688949033
$ws.Range("C4").Value = 8.180485582758843
$ws.Range("D4").Value = 7.770545946825064
$ws.Range("E4").Value = 9.843458967855595
$ws.Range("F4").Value = 40.19174032264964
$ws.Range("I4").Value = 32.12330283803288
$ws.Range("L4").Value = 10.36956841078449
$ws.Range("M4").Value = 16.60650184841725
$ws.Range("B5").Value = 17.8336257496221
$ws.Range("C5").Value = 8.032866124750635
$ws.Range("D5").Value = 7.774436833479853
$ws.Range("E5").Value = 9.839394905585259
$ws.Range("F5").Value = 40.10858002655451
$ws.Range("I5").Value = 32.10264541658203
$ws.Range("L5").Value = 10.3732671441193
$ws.Range("M5").Value = 16.59700625996338
$ws.Range("B6").Value = 17.82000546900234
$ws.Range("C6").Value = 8.008117276074342
$ws.Range("D6").Value = 7.775093877562648
$ws.Range("E6").Value = 9.838715970453334
$ws.Range("F6").Value = 40.09491670895527
$ws.Range("I6").Value = 32.09931419216049
$ws.Range("L6").Value = 10.37390533536542
$ws.Range("M6").Value = 16.59548982568695
$ws.Range("B7").Value = 17.91488223754052
$ws.Range("C7").Value = 8.178510700732824
$ws.Range("D7").Value = 7.770597685174736
$ws.Range("E7").Value = 9.843404433179856
$ws.Range("F7").Value = 40.19060908834178
$ws.Range("I7").Value = 32.12301761548829
$ws.Range("L7").Value = 10.36961668285754
$ws.Range("M7").Value = 16.60636974932477
$ws.Range("B8").Value = 18.34101288855091
$ws.Range("C8").Value = 8.954773425292407
$ws.Range("D8").Value = 7.752355554886138
$ws.Range("E8").Value = 9.863499574267458
$ws.Range("F8").Value = 40.63355198149874
$ws.Range("I8").Value = 32.24196493377384
$ws.Range("L8").Value = 10.35418426784253
$ws.Range("M8").Value = 16.66309708969888
$ws.Range("B9").Value = 19.19510509675683
$ws.Range("C9").Value = 10.35285892843328
$ws.Range("D9").Value = 7.722373731891288
$ws.Range("E9").Value = 9.900791979787304
$ws.Range("F9").Value = 41.57170681545315
$ws.Range("I9").Value = 32.5263790911565
$ws.Range("L9").Value = 10.33648063554029
$ws.Range("M9").Value = 16.80512984684251
$ws.Range("B10").Value = 19.82651826865618
$ws.Range("C10").Value = 11.26912385133704
$ws.Range("D10").Value = 7.70391822810739
$ws.Range("E10").Value = 9.926923476377358
$ws.Range("F10").Value = 42.29871886166547
$ws.Range("I10").Value = 32.76601570514847
$ws.Range("L10").Value = 10.33113318117312
$ws.Range("M10").Value = 16.92768607581569
$ws.Range("B11").Value = 20.11312418218356
$ws.Range("C11").Value = 11.66213646970415
$ws.Range("D11").Value = 7.696307075902159
$ws.Range("E11").Value = 9.938540951006596
$ws.Range("F11").Value = 42.63659173468189
$ws.Range("I11").Value = 32.88154896831792
$ws.Range("L11").Value = 10.33036276751787
$ws.Range("M11").Value = 16.98725223177435
$ws.Range("B12").Value = 20.22144635513894
$ws.Range("C12").Value = 11.80755800301246
$ws.Range("D12").Value = 7.693538484526117
$ws.Range("E12").Value = 9.942901724526667
$ws.Range("F12").Value = 42.76547086945171
$ws.Range("I12").Value = 32.92622100209499
$ws.Range("L12").Value = 10.33030979743401
$ws.Range("M12").Value = 17.01034358829594
$ws.Range("B13").Value = 20.19812838030589
$ws.Range("C13").Value = 11.77638996271164
$ws.Range("D13").Value = 7.694129685637853
$ws.Range("E13").Value = 9.941964262030123
$ws.Range("F13").Value = 42.73767442036437
$ws.Range("I13").Value = 32.91655933657304
$ws.Range("L13").Value = 10.33031059108091
$ws.Range("M13").Value = 17.00534688521634
$ws.Range("B14").Value = 20.1220407686592
$ws.Range("C14").Value = 11.67416849751437
$ws.Range("D14").Value = 7.696077021184553
$ws.Range("E14").Value = 9.938900484602005
$ws.Range("F14").Value = 42.64717652579214
$ws.Range("I14").Value = 32.88520577879117
$ws.Range("L14").Value = 10.33035362609286
$ws.Range("M14").Value = 16.98914133138905
$ws.Range("B15").Value = 20.07540421824613
$ws.Range("C15").Value = 11.61111217471698
$ws.Range("D15").Value = 7.697284636567082
$ws.Range("E15").Value = 9.937018825885833
$ws.Range("F15").Value = 42.5918627601673
$ws.Range("I15").Value = 32.86612041048971
$ws.Range("L15").Value = 10.33041107214059
$ws.Range("M15").Value = 16.97928420926312
$ws.Range("B16").Value = 19.80776497278775
$ws.Range("C16").Value = 11.2429622211188
$ws.Range("D16").Value = 7.704431486034763
$ws.Range("E16").Value = 9.926158841743192
$ws.Range("F16").Value = 42.27677395295827
$ws.Range("I16").Value = 32.75859525900302
$ws.Range("L16").Value = 10.33121695732669
$ws.Range("M16").Value = 16.92386890576586
$ws.Range("B17").Value = 19.6433298806652
$ws.Range("C17").Value = 11.01103153294043
$ws.Range("D17").Value = 7.709017303638932
$ws.Range("E17").Value = 9.919427681537965
$ws.Range("F17").Value = 42.08524226061944
$ws.Range("I17").Value = 32.69429171160528
$ws.Range("L17").Value = 10.33213691385767
$ws.Range("M17").Value = 16.89084104558535
$ws.Range("B18").Value = 19.5486989324862
$ws.Range("C18").Value = 10.87538977111624
$ws.Range("D18").Value = 7.711728723244892
$ws.Range("E18").Value = 9.915530676549631
$ws.Range("F18").Value = 41.97575720334503
$ws.Range("I18").Value = 32.65792087817589
$ws.Range("L18").Value = 10.33282252032165
$ws.Range("M18").Value = 16.87220422350605
$ws.Range("B19").Value = 19.51665326715863
$ws.Range("C19").Value = 10.82907827235356
$ws.Range("D19").Value = 7.712659410628548
$ws.Range("E19").Value = 9.914206842814444
$ws.Range("F19").Value = 41.93880689339179
$ws.Range("I19").Value = 32.64571239329634
$ws.Range("L19").Value = 10.33308153521749
$ws.Range("M19").Value = 16.86595633441655
$ws.Range("B20").Value = 19.66084056312132
$ws.Range("C20").Value = 11.03595278722994
$ws.Range("D20").Value = 7.708521493956932
$ws.Range("E20").Value = 9.920146853737068
$ws.Range("F20").Value = 42.10556155024031
$ws.Range("I20").Value = 32.70107340602817
$ws.Range("L20").Value = 10.33202279020868
$ws.Range("M20").Value = 16.89431976219256
$ws.Range("B21").Value = 20.14439614516365
$ws.Range("C21").Value = 11.70428564972222
$ws.Range("D21").Value = 7.695501952427983
$ws.Range("E21").Value = 9.939801433137632
$ws.Range("F21").Value = 42.67373336214681
$ws.Range("I21").Value = 32.8943901862451
$ws.Range("L21").Value = 10.33033450796929
$ws.Range("M21").Value = 16.99388688970307
$ws.Range("B22").Value = 20.4591567497187
$ws.Range("C22").Value = 12.12124836590292
$ws.Range("D22").Value = 7.687655374302807
$ws.Range("E22").Value = 9.952422498687667
$ws.Range("F22").Value = 43.05046587284382
$ws.Range("I22").Value = 33.02610015720046
$ws.Range("L22").Value = 10.33062266783916
$ws.Range("M22").Value = 17.06207094498916
$ws.Range("B23").Value = 20.29131724473464
$ws.Range("C23").Value = 11.90051594257459
$ws.Range("D23").Value = 7.691782361797331
$ws.Range("E23").Value = 9.945706827118828
$ws.Range("F23").Value = 42.84893414643774
$ws.Range("I23").Value = 32.95531864238356
$ws.Range("L23").Value = 10.33034164957212
$ws.Range("M23").Value = 17.02539986773143
$ws.Range("B24").Value = 19.65292427046247
$ws.Range("C24").Value = 11.02469305892787
$ws.Range("D24").Value = 7.708745416053396
$ws.Range("E24").Value = 9.919821800640527
$ws.Range("F24").Value = 42.09637323287214
$ws.Range("I24").Value = 32.69800553776125
$ws.Range("L24").Value = 10.33207389737198
$ws.Range("M24").Value = 16.89274593906384
$ws.Range("B25").Value = 18.96284380967248
$ws.Range("C25").Value = 9.994205526450182
$ws.Range("D25").Value = 7.729860346167293
$ws.Range("E25").Value = 9.890927500743702
$ws.Range("F25").Value = 41.31093064449718
$ws.Range("I25").Value = 32.44400619024821
$ws.Range("L25").Value = 10.33992475564365
$ws.Range("M25").Value = 16.76346434626922
